$d = $word.ActiveDocument

# Fix 1: "see snail" -> "sea snail"
$d.Content.Find.Execute("An abalone is a type of see snail. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "An abalone is a type of sea snail. ", 2)

# Fix 2: "For the dataset we are using the real values are known (i." -> "In the dataset we are using the true number of rings is known (i."
$d.Content.Find.Execute("For the dataset we are using the real values are known (i.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "In the dataset we are using the true number of rings is known (i.", 2)
